$p = $ppt.ActivePresentation

# Delete the four slides (originally sldId 276, 277, 278, 279) that covered the
# Flask "used car price predictor" implementation details. They sit at
# positions 8-11 (1-based) in the original slide order. Deleting from the
# highest index down keeps the remaining indices stable while we work.
$p.Slides.Item(11).Delete()
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()
$p.Slides.Item(8).Delete()
